# "excel file for login function is added"
# Fill the Status column (C) on Sheet1 with Pass/Fail results:
# all sign-up rows fail the login check except the last one (row 17),
# which passes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 16; $r++) {
    $ws.Range("C" + $r).Value = "Fail"
}
$ws.Range("C17").Value = "Pass"
